# Update cryptocurrency price / 1h-volume-change table with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.184.33"
$ws.Range("E2").Value = "  -4.91%  "

$ws.Range("D3").Value = "3.316.47"
$ws.Range("E3").Value = "  -5.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.98"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.37%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.317.03"
$ws.Range("E8").Value = "  -5.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.477"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.38%  "

$ws.Range("E10").Value = "  -5.33%  "

$ws.Range("E11").Value = "  -5.53%  "

$ws.Range("E12").Value = "  -4.20%  "

$ws.Range("D13").Value = "3.880.84"
$ws.Range("E13").Value = "  -5.18%  "

$ws.Range("E14").Value = "  -1.30%  "

$ws.Range("D15").Value = "3.315.69"
$ws.Range("E15").Value = "  -5.18%  "

$ws.Range("E16").Value = "  -6.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.83"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").Value = "61.256.17"
$ws.Range("E18").Value = "  -4.80%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.23%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.09"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.18"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "353.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -8.88%  "

$ws.Range("E23").Value = "  -4.56%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").Value = "3.447.02"
$ws.Range("E25").Value = "  -5.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.14"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.55%  "

$ws.Range("E27").Value = "  -7.11%  "

$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.19"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.48"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("E31").Value = "  -3.25%  "

$ws.Range("E32").Value = "  -6.65%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("E34").Value = "  -4.92%  "

$ws.Range("D35").Value = "3.341.99"
$ws.Range("E35").Value = "  -5.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.45"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.34"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.35%  "

$ws.Range("E38").Value = "  -2.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.94"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0755"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.01%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.85"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.86%  "

$ws.Range("E44").Value = "  -7.47%  "

$ws.Range("E45").Value = "  -4.20%  "

$ws.Range("E46").Value = "  -5.91%  "

$ws.Range("E47").Value = "  -6.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.29"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -10.33%  "

$ws.Range("E49").Value = "  -3.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.854"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.43%  "

$ws.Range("D51").Value = "2.196.24"
$ws.Range("E51").Value = "  -8.57%  "

